# Applies the "Elimina antiguos EC y agrega nuevos y modifica Antigua BD" update:
#   - Updates the overdue-balance total (E11) and the period count (F13).
#   - Fixes the "Novedad de Retiro" / "Novedad de Ingreso" column headers (H15/I15 were swapped).
#   - Adds a new account-statement line for period 2509 (copy of the 2508 row),
#     inserted right below the existing last data row, and re-flows the row
#     that used to be "last" back to a normal interior row.
#   - Inserts a new row before the signature block so the underline and the
#     "NOMBRE DEL REPRESENTANTE LEGAL" caption each get their own row again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) Header totals -----------------------------------------------------
$ws.Range("E11").Value = 412884
$ws.Range("F13").Value = 6

# ---- 2) Swap the Retiro/Ingreso column headers -----------------------------
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

# ---- 3) Add the new 2509 period row ----------------------------------------
# Push everything from row 21 down by one, leaving a blank row 21 to build the
# new statement line in (this is also what shifts the old rows 25/26 down to
# 26/27 automatically).
$ws.Rows.Item(21).Insert()

# Clone row 20's full formatting + values into the freshly inserted row 21 --
# row 20 currently holds the most-recent (2508) statement line, which is
# exactly the template the new 2509 line should look like.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$ws.Range("B21:J21").Value = $ws.Range("B20:J20").Value()

# The new row represents period 2509 instead of 2508.
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2509"

# ---- 4) Row 20 is no longer the last row, so it goes back to the regular
#         interior-row formatting (cloning row 19's look, keeping its values).
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Column E (the period code) is centered in the refreshed layout.
$ws.Range("E16:E21").HorizontalAlignment = -4108

$ws.Application.CutCopyMode = 0
